$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing
# Usain/Ishan/Nehal/Papa table from A:D to B:E (and preserves the
# explicit-black-font style that lived on the old "Nehal" column,
# which lands on the new "Sham" column at the same relative offset).
$ws.Columns.Item(1).Insert()

# New column A: row labels.
$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Race"
$ws.Range("A3").Value = "Time"
$ws.Range("A4").Value = "Last weight"
$ws.Range("A5").Value = "Current weight"
$ws.Range("A6").Value = "Past BMI"
$ws.Range("A7").Value = "Current BMI"

# Relabel the header row for the shifted columns.
$ws.Range("D1").Value = "Sham"
$ws.Range("E1").Value = "Ram"

# New column F: header + data (rows 2-5 mirror column E's values).
$ws.Range("F1").Value = "Laxman"
$ws.Range("F2").Value = 100
$ws.Range("F3").Value = 29
$ws.Range("F4").Value = 171
$ws.Range("F5").Value = 172

# New rows 6 and 7.
$ws.Range("B6").Value = 24
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = 29
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20

$ws.Range("B7").Value = 24
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = 32
$ws.Range("E7").Value = 27
$ws.Range("F7").Value = 20

# Restore the selection/active cell as saved in the source file.
$ws.Range("H9").Select()
